$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8093
$ws1.Range("F5").Value = 5908
$ws1.Range("F7").Value = 92
$ws1.Range("F11").Value = 431

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8093
$ws4.Range("F5").Value = 5908
$ws4.Range("F7").Value = 92
$ws4.Range("F15").Value = 431
